$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the original cell/value that is no longer present in the target workbook
$ws.Range("G7").ClearContents()

# Set cell values in the exact order needed so the shared-string table
# is rebuilt with the same index ordering as the target file.
$ws.Range("I8").Value  = "sdfsdf"
$ws.Range("G8").Value  = "sdflkjsdlfj"
$ws.Range("G4").Value  = "sldkfjsdlfkdsjf"
$ws.Range("J4").Value  = "sldkjf"
$ws.Range("K10").Value = "sdlfjslfkj"
$ws.Range("K6").Value  = "sdlkfjslfkj"
$ws.Range("I5").Value  = "sdlkfsdlfkj"
$ws.Range("I2").Value  = "sldkfj"
$ws.Range("H6").Value  = "sldkfdslfkj"
$ws.Range("K8").Value  = "sdlkfjsdlkfj"
$ws.Range("I10").Value = "slkfdslfkj"
$ws.Range("L2").Value  = "sdlkfsdlkfj"
$ws.Range("L6").Value  = "sldkjf"

# Move selection to L2 to match the target view state
$ws.Range("L2").Select()
